$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing existing rows 4-9 down to 5-10.
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the weekly data entry.
$ws.Range("A4").Value = 8
$ws.Range("B4").Value = "Terminal La Palmera de La Serena"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 44953
$ws.Range("D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 100112039
$ws.Range("G4").Value = "Ciboulette"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 2500
$ws.Range("M4").Value = 2250
$ws.Range("N4").Value = "`$/docena de atados"
$ws.Range("O4").Value = "Provincia del Elquí"
$ws.Range("P4").Value = 750
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = "Hortaliza"
